$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 7
$ws.Range("K7").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB_MAIN"
$ws.Range("L7").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB"
$ws.Range("N7").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_MAIN"

# Row 8
$ws.Range("K8").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB_MAIN"
$ws.Range("L8").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB"
$ws.Range("N8").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_SUB_MAIN"
$ws.Range("O8").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_SUB"

# Row 10
$ws.Range("K10").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB_MAIN"
$ws.Range("L10").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB"

# Row 11
$ws.Range("N11").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_MAIN"

# Row 12
$ws.Range("K12").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB_MAIN"
$ws.Range("L12").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB"
$ws.Range("N12").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_MAIN"

# Row 13
$ws.Range("K13").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB_MAIN"
$ws.Range("L13").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB"
$ws.Range("N13").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_SUB_MAIN"
$ws.Range("O13").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_SUB"

# Row 15
$ws.Range("K15").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB_MAIN"
$ws.Range("L15").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB"

# Row 16
$ws.Range("N16").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_MAIN"

# Row 17
$ws.Range("K17").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB_MAIN"
$ws.Range("L17").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB"
$ws.Range("N17").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_MAIN"

# Row 18
$ws.Range("K18").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB_MAIN"
$ws.Range("L18").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB"
$ws.Range("N18").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_SUB_MAIN"
$ws.Range("O18").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_SUB"

# Row 20
$ws.Range("K20").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB_MAIN"
$ws.Range("L20").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB"

# Row 21
$ws.Range("N21").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_MAIN"

# Row 22
$ws.Range("K22").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB_MAIN"
$ws.Range("L22").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB"
$ws.Range("N22").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_MAIN"

# Row 23
$ws.Range("K23").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB_MAIN"
$ws.Range("L23").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB"
$ws.Range("N23").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_SUB_MAIN"
$ws.Range("O23").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_SUB"

# Row 25
$ws.Range("K25").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB_MAIN"
$ws.Range("L25").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB"

# Row 26
$ws.Range("N26").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_MAIN"

# Row 27
$ws.Range("K27").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB_MAIN"
$ws.Range("L27").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB"
$ws.Range("N27").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_MAIN"

# Row 28
$ws.Range("K28").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB_MAIN"
$ws.Range("L28").Value = "MPA_TRANSFER_SCENARIO_SENDER_SUB"
$ws.Range("N28").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_SUB_MAIN"
$ws.Range("O28").Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_SUB"
